$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 185.90909
$ws.Range("I5").Value = 185.90909
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 185.90909
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -70.90908999999999
$ws.Range("N5").ClearContents()

$ws.Range("H88").Value = 1097.5
$ws.Range("I88").Value = 255.16667
$ws.Range("K88").Value = 255.16667
$ws.Range("M88").Value = 150.83333

$ws.Range("H91").Value = 1097.5
$ws.Range("I91").Value = 255.16667
$ws.Range("K91").Value = 255.16667
$ws.Range("M91").Value = 1148.83333

$ws.Range("H112").Value = 1999.4445
$ws.Range("J112").Value = 1999.4445
$ws.Range("L112").Value = 5998.333500000001
$ws.Range("N112").Value = -8214.333500000001

$ws.Range("H115").Value = 818.4545000000001
$ws.Range("I115").Value = 722
$ws.Range("K115").Value = 2166
$ws.Range("M115").Value = -599

$ws.Range("H116").Value = 22891.273
$ws.Range("I116").Value = 16974.875
$ws.Range("K116").Value = 16974.875
$ws.Range("M116").Value = -13532.875

$ws.Range("H132").Value = 2293.1667
$ws.Range("I132").Value = 2314.0637
$ws.Range("J132").Value = 2152.8572
$ws.Range("K132").Value = 6942.1911
$ws.Range("L132").Value = 6458.571599999999
$ws.Range("M132").Value = -4412.1911
$ws.Range("N132").Value = -11518.5716

$ws.Range("H137").Value = 6210.9473
$ws.Range("J137").Value = 9586.857
$ws.Range("L137").Value = 28760.571
$ws.Range("N137").Value = -33860.571

$ws.Range("H138").Value = 6747.85
$ws.Range("J138").Value = 6942.4155
$ws.Range("L138").Value = 20827.2465
$ws.Range("N138").Value = -31107.2465

$ws.Range("H141").Value = 13797.25
$ws.Range("I141").Value = 13797.25
$ws.Range("K141").Value = 41391.75
$ws.Range("M141").Value = -36211.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7466.1
$ws.Range("I32").Value = 5350.109
$ws.Range("K32").Value = 5350.109
$ws.Range("M32").Value = -5063.109

$ws.Range("H43").Value = 17871.125
$ws.Range("J43").Value = 20354.5
$ws.Range("L43").Value = 20354.5
$ws.Range("N43").Value = -20980.5

$ws.Range("H109").Value = 132597.83
$ws.Range("J109").Value = 132597.83
$ws.Range("L109").Value = 132597.83
$ws.Range("N109").Value = -135371.83

$ws.Range("H122").Value = 3657.2354
$ws.Range("I122").Value = 2521.75
$ws.Range("J122").Value = 4666.5557
$ws.Range("K122").Value = 7565.25
$ws.Range("L122").Value = 13999.6671
$ws.Range("M122").Value = -5115.25
$ws.Range("N122").Value = -18899.6671

$ws.Range("H132").Value = 5973.488
$ws.Range("I132").Value = 5100
$ws.Range("J132").Value = 7487.533
$ws.Range("K132").Value = 15300
$ws.Range("L132").Value = 22462.599
$ws.Range("M132").Value = -12770
$ws.Range("N132").Value = -27522.599

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4801.0938
$ws.Range("I107").Value = 5472.8076
$ws.Range("K107").Value = 5472.8076
$ws.Range("M107").Value = -3552.8076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 24369.545
$ws.Range("I16").Value = 8912.857
$ws.Range("K16").Value = 8912.857
$ws.Range("M16").Value = -8625.857

$ws.Range("H113").Value = 24369.545
$ws.Range("I113").Value = 8912.857
$ws.Range("K113").Value = 8912.857
$ws.Range("M113").Value = -6742.857

$ws.Range("H132").Value = 4384.8
$ws.Range("I132").Value = 3474.6667
$ws.Range("J132").Value = 5750
$ws.Range("K132").Value = 10424.0001
$ws.Range("L132").Value = 17250
$ws.Range("M132").Value = -7894.000100000001
$ws.Range("N132").Value = -22310

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 123636.18
$ws.Range("I131").Value = 170881.5
$ws.Range("K131").Value = 512644.5
$ws.Range("M131").Value = -507604.5

$ws.Range("H132").Value = 920627.9399999999
$ws.Range("I132").Value = 168806.67
$ws.Range("K132").Value = 1519260.03
$ws.Range("M132").Value = -1516730.03

$ws.Range("H140").Value = 5444.8076
$ws.Range("I140").Value = 3674.5386
$ws.Range("K140").Value = 11023.6158
$ws.Range("M140").Value = -5843.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 38420.75
$ws.Range("J95").Value = 38420.75
$ws.Range("L95").Value = 38420.75
$ws.Range("N95").Value = -43912.75

$ws.Range("H102").Value = 2586.525
$ws.Range("I102").Value = 1583.8334
$ws.Range("K102").Value = 1583.8334
$ws.Range("M102").Value = 38.16660000000002

$ws.Range("H113").Value = 1115390
$ws.Range("I113").Value = 3336336.8
$ws.Range("J113").Value = 4916.6665
$ws.Range("K113").Value = 3336336.8
$ws.Range("L113").Value = 4916.6665
$ws.Range("M113").Value = -3334166.8
$ws.Range("N113").Value = -9256.666499999999

$ws.Range("H122").Value = 3528.889
$ws.Range("I122").Value = 2717.2
$ws.Range("J122").Value = 4543.5
$ws.Range("K122").Value = 8151.599999999999
$ws.Range("L122").Value = 13630.5
$ws.Range("M122").Value = -5701.599999999999
$ws.Range("N122").Value = -18530.5

$ws.Range("H125").Value = 90000
$ws.Range("J125").Value = 90000
$ws.Range("L125").Value = 90000
$ws.Range("N125").Value = -94920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4454.1113
$ws.Range("I40").Value = 4454.1113
$ws.Range("K40").Value = 4454.1113
$ws.Range("M40").Value = -4318.1113

$ws.Range("H61").Value = 18422
$ws.Range("I61").Value = 18866.438
$ws.Range("K61").Value = 18866.438
$ws.Range("M61").Value = -18664.438

$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872

$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360

$ws.Range("H113").Value = 18422
$ws.Range("I113").Value = 18866.438
$ws.Range("K113").Value = 18866.438
$ws.Range("M113").Value = -16696.438

$ws.Range("H124").Value = 99994
$ws.Range("J124").Value = 99994
$ws.Range("L124").Value = 99994
$ws.Range("N124").Value = -109814

$ws.Range("H132").Value = 6676.3335
$ws.Range("I132").Value = 6155.2856
$ws.Range("K132").Value = 18465.8568
$ws.Range("M132").Value = -15935.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 3004.2
$ws.Range("J48").Value = 3004.2
$ws.Range("L48").Value = 3004.2
$ws.Range("N48").Value = -4142.2

$ws.Range("H86").Value = 92616.664
$ws.Range("J86").Value = 92616.664
$ws.Range("L86").Value = 92616.664
$ws.Range("N86").Value = -94862.664

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H89").Value = 92616.664
$ws.Range("J89").Value = 92616.664
$ws.Range("L89").Value = 463083.32
$ws.Range("N89").Value = -474315.32

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H96").Value = 57134.832
$ws.Range("I96").Value = 92439.17999999999
$ws.Range("J96").Value = 1656.5714
$ws.Range("K96").Value = 92439.17999999999
$ws.Range("L96").Value = 1656.5714
$ws.Range("M96").Value = -91066.17999999999
$ws.Range("N96").Value = -4402.5714

$ws.Range("H125").Value = 68664
$ws.Range("J125").Value = 68664
$ws.Range("L125").Value = 68664
$ws.Range("N125").Value = -78504

$ws.Range("H132").Value = 36151.234
$ws.Range("I132").Value = 2211.28
$ws.Range("J132").Value = 205851
$ws.Range("K132").Value = 6633.84
$ws.Range("L132").Value = 617553
$ws.Range("M132").Value = -4103.84
$ws.Range("N132").Value = -622613

$ws.Range("H135").Value = 102499.5
$ws.Range("J135").Value = 102499.5
$ws.Range("L135").Value = 102499.5
$ws.Range("N135").Value = -112639.5
